$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 3 (item id 18511)
$ws.Range("H3").Value = 65000
$ws.Range("J3").Value = 65000
$ws.Range("L3").Value = 65000
$ws.Range("N3").Value = -65228
# Row 38 (item id 4599)
$ws.Range("H38").Value = 1994.7273
$ws.Range("I38").Value = 157.5
$ws.Range("K38").Value = 472.5
$ws.Range("M38").Value = -100.5
# Row 51 (item id 5486)
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()
# Row 70 (item id 12604)
$ws.Range("H70").Value = 1749.6666
$ws.Range("I70").Value = 1749
$ws.Range("J70").Value = 1750
$ws.Range("K70").Value = 5247
$ws.Range("L70").Value = 5250
$ws.Range("M70").Value = -4977
$ws.Range("N70").Value = -5790
# Row 73 (item id 12604)
$ws.Range("H73").Value = 1749.6666
$ws.Range("I73").Value = 1749
$ws.Range("J73").Value = 1750
$ws.Range("K73").Value = 5247
$ws.Range("L73").Value = 5250
$ws.Range("M73").Value = -4311
$ws.Range("N73").Value = -7122
# Row 74 (item id 5507)
$ws.Range("H74").Value = 0
$ws.Range("I74").Value = 0
$ws.Range("K74").Value = 0
$ws.Range("M74").ClearContents()
# Row 77 (item id 5507)
$ws.Range("H77").Value = 0
$ws.Range("I77").Value = 0
$ws.Range("K77").Value = 0
$ws.Range("M77").ClearContents()
# Row 102 (item id 18511)
$ws.Range("H102").Value = 65000
$ws.Range("J102").Value = 65000
$ws.Range("L102").Value = 65000
$ws.Range("N102").Value = -71490
# Row 138 (item id 44169)
$ws.Range("H138").Value = 1598.5

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 61 (item id 43999)
$ws.Range("H61").Value = 3124
$ws.Range("I61").Value = 3278
$ws.Range("J61").Value = 2200
$ws.Range("K61").Value = 3278
$ws.Range("L61").Value = 2200
$ws.Range("M61").Value = -3066
$ws.Range("N61").Value = -2624
# Row 74 (item id 44000)
$ws.Range("H74").Value = 1862.25
$ws.Range("I74").Value = 1414
$ws.Range("K74").Value = 1414
$ws.Range("M74").Value = -540
# Row 77 (item id 44000)
$ws.Range("H77").Value = 1862.25
$ws.Range("I77").Value = 1414
$ws.Range("K77").Value = 7070
$ws.Range("M77").Value = -2702
# Row 102 (item id 19945)
$ws.Range("H102").Value = 5192
$ws.Range("I102").Value = 5192
$ws.Range("K102").Value = 5192
$ws.Range("M102").Value = -3570
# Row 136 (item id 43999)
$ws.Range("H136").Value = 3124
$ws.Range("I136").Value = 3278
$ws.Range("J136").Value = 2200
$ws.Range("K136").Value = 9834
$ws.Range("L136").Value = 6600
$ws.Range("M136").Value = -7284
$ws.Range("N136").Value = -11700

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 75 (item id 11872)
$ws.Range("H75").Value = 5310.5557
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
# Row 78 (item id 11872)
$ws.Range("H78").Value = 5310.5557
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
# Row 99 (item id 19943)
$ws.Range("H99").Value = 1491
$ws.Range("I99").Value = 1491
$ws.Range("K99").Value = 1491
$ws.Range("M99").Value = 7

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 16 (item id 27691)
$ws.Range("H16").Value = 881
$ws.Range("I16").Value = 881
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 881
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -594
$ws.Range("N16").ClearContents()
# Row 58 (item id 44021)
$ws.Range("H58").Value = 3287.4285
$ws.Range("I58").Value = 3287.4285
$ws.Range("K58").Value = 3287.4285
$ws.Range("M58").Value = -3084.4285
# Row 113 (item id 27691)
$ws.Range("H113").Value = 881
$ws.Range("I113").Value = 881
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 881
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1289
$ws.Range("N113").ClearContents()
# Row 132 (item id 44019)
$ws.Range("H132").Value = 3464
$ws.Range("I132").Value = 2511.3333
$ws.Range("J132").Value = 4416.6665
$ws.Range("K132").Value = 7533.999899999999
$ws.Range("L132").Value = 13249.9995
$ws.Range("M132").Value = -5003.999899999999
$ws.Range("N132").Value = -18309.9995
# Row 134 (item id 44020)
$ws.Range("H134").Value = 1466.3334
$ws.Range("J134").Value = 1199
$ws.Range("L134").Value = 3597
$ws.Range("N134").Value = -8667
# Row 136 (item id 44021)
$ws.Range("H136").Value = 3287.4285
$ws.Range("I136").Value = 3287.4285
$ws.Range("K136").Value = 9862.2855
$ws.Range("M136").Value = -7312.2855

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 2 (item id 5062)
$ws.Range("H2").Value = 7.2
$ws.Range("I2").Value = 7.2
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 7.2
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 105.8
$ws.Range("N2").ClearContents()
# Row 3 (item id 4091)
$ws.Range("H3").Value = 5541833.5
$ws.Range("I3").Value = 5812750
$ws.Range("K3").Value = 5812750
$ws.Range("M3").Value = -5812634
# Row 25 (item id 2687)
$ws.Range("H25").Value = 1250
$ws.Range("J25").Value = 1250
$ws.Range("L25").Value = 1250
$ws.Range("N25").Value = -2308
# Row 80 (item id 12521)
$ws.Range("H80").Value = 20000
$ws.Range("I80").Value = 20000
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 20000
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -19002
$ws.Range("N80").ClearContents()
# Row 83 (item id 12521)
$ws.Range("H83").Value = 20000
$ws.Range("I83").Value = 20000
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 100000
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -95008
$ws.Range("N83").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 42 (item id 4333)
$ws.Range("H42").Value = 25024
$ws.Range("I42").Value = 25024
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 25024
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -24461
$ws.Range("N42").ClearContents()
# Row 46 (item id 5282)
$ws.Range("H46").Value = 490
$ws.Range("J46").Value = 490
$ws.Range("L46").Value = 490
$ws.Range("N46").Value = -866
# Row 49 (item id 4333)
$ws.Range("H49").Value = 25024
$ws.Range("I49").Value = 25024
$ws.Range("J49").Value = 0
$ws.Range("K49").Value = 25024
$ws.Range("L49").Value = 0
$ws.Range("M49").Value = -24877
$ws.Range("N49").ClearContents()
# Row 50 (item id 3426)
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()
# Row 82 (item id 12565)
$ws.Range("H82").Value = 0
$ws.Range("I82").Value = 0
$ws.Range("K82").Value = 0
$ws.Range("M82").ClearContents()
# Row 85 (item id 12565)
$ws.Range("H85").Value = 0
$ws.Range("I85").Value = 0
$ws.Range("K85").Value = 0
$ws.Range("M85").ClearContents()
# Row 104 (item id 18675)
$ws.Range("H104").Value = 31435
$ws.Range("J104").Value = 31435
$ws.Range("L104").Value = 31435
$ws.Range("N104").Value = -38423
# Row 122 (item id 36247)
$ws.Range("H122").Value = 4739.533
$ws.Range("I122").Value = 4739.533
$ws.Range("K122").Value = 14218.599
$ws.Range("M122").Value = -11768.599

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 122 (item id 36208)
$ws.Range("H122").Value = 1665
$ws.Range("I122").Value = 1665
$ws.Range("K122").Value = 4995
$ws.Range("M122").Value = -2545
# Row 126 (item id 36210)
$ws.Range("H126").Value = 1475.3334
$ws.Range("I126").Value = 1475.3334
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 4426.0002
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -1956.0002
$ws.Range("N126").ClearContents()

